# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-data refresh to the Excalibur_Profits workbook:
# updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the affected
# leve rows across the ALC, ARM, BSM, CRP, CUL, GSM and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 147.625
$ws.Range("I6").Value = 97.42856999999999
$ws.Range("J6").Value = 499
$ws.Range("K6").Value = 292.28571
$ws.Range("L6").Value = 1497
$ws.Range("M6").Value = -180.28571
$ws.Range("N6").Value = -1721

$ws.Range("H64").Value = 7821.923
$ws.Range("I64").Value = 4625
$ws.Range("J64").Value = 9242.777
$ws.Range("K64").Value = 4625
$ws.Range("L64").Value = 9242.777
$ws.Range("M64").Value = -4377
$ws.Range("N64").Value = -9738.777

$ws.Range("H67").Value = 7821.923
$ws.Range("I67").Value = 4625
$ws.Range("J67").Value = 9242.777
$ws.Range("K67").Value = 4625
$ws.Range("L67").Value = 9242.777
$ws.Range("M67").Value = -3767
$ws.Range("N67").Value = -10958.777

$ws.Range("H69").Value = 7146.8276
$ws.Range("I69").Value = 4323.3076
$ws.Range("K69").Value = 12969.9228
$ws.Range("M69").Value = -12095.9228

$ws.Range("H72").Value = 7146.8276
$ws.Range("I72").Value = 4323.3076
$ws.Range("K72").Value = 38909.7684
$ws.Range("M72").Value = -34541.7684

$ws.Range("H111").Value = 464.5
$ws.Range("I111").Value = 464.5
$ws.Range("K111").Value = 1393.5
$ws.Range("M111").Value = 1673.5

$ws.Range("H129").Value = 2024.9445
$ws.Range("I129").Value = 1396.4286
$ws.Range("K129").Value = 4189.2858
$ws.Range("M129").Value = 810.7142000000003

$ws.Range("H131").Value = 4069.8518
$ws.Range("I131").Value = 3375.4119
$ws.Range("J131").Value = 5250.4
$ws.Range("K131").Value = 10126.2357
$ws.Range("L131").Value = 15751.2
$ws.Range("M131").Value = -5086.235700000001
$ws.Range("N131").Value = -25831.2

$ws.Range("H132").Value = 2659
$ws.Range("I132").Value = 2605.4
$ws.Range("K132").Value = 7816.200000000001
$ws.Range("M132").Value = -5286.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5854.433
$ws.Range("I32").Value = 6058.3213
$ws.Range("K32").Value = 6058.3213
$ws.Range("M32").Value = -5771.3213

$ws.Range("H61").Value = 3877.4055
$ws.Range("I61").Value = 2427.2
$ws.Range("J61").Value = 6898.6665
$ws.Range("K61").Value = 2427.2
$ws.Range("L61").Value = 6898.6665
$ws.Range("M61").Value = -2215.2
$ws.Range("N61").Value = -7322.6665

$ws.Range("H74").Value = 2390.1082
$ws.Range("I74").Value = 1600.4517
$ws.Range("J74").Value = 6470
$ws.Range("K74").Value = 1600.4517
$ws.Range("L74").Value = 6470
$ws.Range("M74").Value = -726.4517000000001
$ws.Range("N74").Value = -8218

$ws.Range("H77").Value = 2390.1082
$ws.Range("I77").Value = 1600.4517
$ws.Range("J77").Value = 6470
$ws.Range("K77").Value = 8002.2585
$ws.Range("L77").Value = 32350
$ws.Range("M77").Value = -3634.2585
$ws.Range("N77").Value = -41086

$ws.Range("H122").Value = 4416.0527
$ws.Range("I122").Value = 4241.875
$ws.Range("K122").Value = 12725.625
$ws.Range("M122").Value = -10275.625

$ws.Range("H132").Value = 3692.7273
$ws.Range("I132").Value = 1981.193
$ws.Range("K132").Value = 5943.579
$ws.Range("M132").Value = -3413.579

$ws.Range("H136").Value = 3877.4055
$ws.Range("I136").Value = 2427.2
$ws.Range("J136").Value = 6898.6665
$ws.Range("K136").Value = 7281.599999999999
$ws.Range("L136").Value = 20695.9995
$ws.Range("M136").Value = -4731.599999999999
$ws.Range("N136").Value = -25795.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2626.5
$ws.Range("I20").Value = 3212.8
$ws.Range("J20").Value = 2207.7144
$ws.Range("K20").Value = 3212.8
$ws.Range("L20").Value = 2207.7144
$ws.Range("M20").Value = -2965.8
$ws.Range("N20").Value = -2701.7144

$ws.Range("H100").Value = 18303
$ws.Range("J100").Value = 18303
$ws.Range("L100").Value = 18303
$ws.Range("N100").Value = -20467

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H31").Value = 4370.185
$ws.Range("I31").Value = 2957.3462
$ws.Range("J31").Value = 5682.107
$ws.Range("K31").Value = 2957.3462
$ws.Range("L31").Value = 5682.107
$ws.Range("M31").Value = -2662.3462
$ws.Range("N31").Value = -6272.107

$ws.Range("H34").Value = 4370.185
$ws.Range("I34").Value = 2957.3462
$ws.Range("J34").Value = 5682.107
$ws.Range("K34").Value = 2957.3462
$ws.Range("L34").Value = 5682.107
$ws.Range("M34").Value = -2755.3462
$ws.Range("N34").Value = -6086.107

$ws.Range("H43").Value = 100339.164
$ws.Range("J43").Value = 100339.164
$ws.Range("L43").Value = 100339.164
$ws.Range("N43").Value = -100707.164

$ws.Range("H101").Value = 100339.164
$ws.Range("J101").Value = 100339.164
$ws.Range("L101").Value = 100339.164
$ws.Range("N101").Value = -106829.164

$ws.Range("H132").Value = 1618.25
$ws.Range("I132").Value = 1595.05
$ws.Range("J132").Value = 1734.25
$ws.Range("K132").Value = 4785.15
$ws.Range("L132").Value = 5202.75
$ws.Range("M132").Value = -2255.15
$ws.Range("N132").Value = -10262.75

$ws.Range("H134").Value = 1944.8889
$ws.Range("I134").Value = 1937.5217
$ws.Range("K134").Value = 5812.5651
$ws.Range("M134").Value = -3277.5651

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 14.833333
$ws.Range("I6").Value = 14.833333
$ws.Range("K6").Value = 44.499999
$ws.Range("M6").Value = 68.500001

$ws.Range("H7").Value = 3154003.5
$ws.Range("J7").Value = 10500000
$ws.Range("L7").Value = 31500000
$ws.Range("N7").Value = -31500224

$ws.Range("H10").Value = 493.25
$ws.Range("I10").Value = 474.33334
$ws.Range("K10").Value = 1423.00002
$ws.Range("M10").Value = -1284.00002

$ws.Range("H131").Value = 1335.5294
$ws.Range("J131").Value = 2187.25
$ws.Range("L131").Value = 6561.75
$ws.Range("N131").Value = -16641.75

$ws.Range("H133").Value = 4932.857
$ws.Range("I133").Value = 4765
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 14295
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -9235
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H132").Value = 1164.5714
$ws.Range("I132").Value = 1114.9166
$ws.Range("J132").Value = 1462.5
$ws.Range("K132").Value = 3344.7498
$ws.Range("L132").Value = 4387.5
$ws.Range("M132").Value = -814.7498000000001
$ws.Range("N132").Value = -9447.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1032.8182
$ws.Range("I81").Value = 880.2857
$ws.Range("J81").Value = 1299.75
$ws.Range("K81").Value = 1760.5714
$ws.Range("L81").Value = 2599.5
$ws.Range("M81").Value = -699.5714
$ws.Range("N81").Value = -4721.5

$ws.Range("H84").Value = 1032.8182
$ws.Range("I84").Value = 880.2857
$ws.Range("J84").Value = 1299.75
$ws.Range("K84").Value = 8802.857
$ws.Range("L84").Value = 12997.5
$ws.Range("M84").Value = -3498.857
$ws.Range("N84").Value = -23605.5

$ws.Range("H107").Value = 2281.0386
$ws.Range("I107").Value = 1436.0625
$ws.Range("J107").Value = 3633
$ws.Range("K107").Value = 4308.1875
$ws.Range("L107").Value = 10899
$ws.Range("M107").Value = -2388.1875
$ws.Range("N107").Value = -14739

$ws.Range("H132").Value = 15000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 45000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -50060

$ws.Range("H136").Value = 11522013
$ws.Range("I136").Value = 16457690
$ws.Range("J136").Value = 5433
$ws.Range("K136").Value = 49373070
$ws.Range("L136").Value = 16299
$ws.Range("M136").Value = -49370520
$ws.Range("N136").Value = -21399

Write-Output "Applied scheduled market-data refresh."
